# Remove the three "Ροδάκινα® Ναουσας (Ζυγιζόμενο) /Kgr" (peach) line items
# (original rows 3, 4 and 7) from the daily price-list report, then refresh
# the per-row running totals (K/L) and the grand-total row for the rows
# that shift up, matching the regenerated report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so the remaining row numbers don't need to be recomputed
# as we go.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Rows that shifted up now carry refreshed SalesQuantity / Turnover figures.
$ws.Cells.Item(5, 11).Value = 9.865
$ws.Cells.Item(5, 12).Value = 14.4
$ws.Cells.Item(6, 11).Value = 15
$ws.Cells.Item(6, 12).Value = 39.89
$ws.Cells.Item(7, 11).Value = 53.497
$ws.Cells.Item(7, 12).Value = 40.85

# Grand-total row (now row 8 after the three deletions).
$ws.Cells.Item(8, 11).Value = 89.362
$ws.Cells.Item(8, 12).Value = 112.06

# The two colour-scale conditional formats covered I1:I10 / J1:J10; with
# only 7 data rows left they should cover I1:I7 / J1:J7 instead.
$ws.Range("I1:I10").FormatConditions.Delete()
$ws.Range("J1:J10").FormatConditions.Delete()

$null = $ws.Range("I1:I7").FormatConditions.AddColorScale(3)
$null = $ws.Range("J1:J7").FormatConditions.AddColorScale(3)

Write-Host "Removed peach rows, refreshed totals, and resized conditional formatting."
